$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) / Volume(1h) (E) columns with the latest scrape.
# Some Price values (e.g. "1.00", "584.05") look numeric, so a plain
# Value assignment would make Excel coerce them into Number cells and
# drop the significant trailing zero. Prefixing with a leading apostrophe
# (PowerShell single-quoted '' -> literal ') is the standard Excel
# text-entry trick that forces those cells to stay literal Text, exactly
# like the source data.
$ws.Range('D2').Value = '62.950.93'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '2.571.78'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('D5').Value = '''584.05'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').Value = '''143.99'
$ws.Range('E6').Value = '  -2.31%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '''0.589'
$ws.Range('E8').Value = '  -2.06%  '
$ws.Range('E9').Value = '  -2.45%  '
$ws.Range('D10').Value = '''5.57'
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('E12').Value = '  -2.23%  '
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('D14').Value = '3.032.57'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').Value = '62.861.76'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('E16').Value = '  -1.93%  '
$ws.Range('D17').Value = '2.572.04'
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('E18').Value = '  -2.95%  '
$ws.Range('D19').Value = '''340.32'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('D21').Value = '''6.62'
$ws.Range('E21').Value = '  -3.80%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('E23').Value = '  +3.40%  '
$ws.Range('D24').Value = '''67.61'
$ws.Range('E24').Value = '  +1.11%  '
$ws.Range('E25').Value = '  +7.68%  '
$ws.Range('E26').Value = '  -2.48%  '
$ws.Range('D27').Value = '''0.165'
$ws.Range('E27').Value = '  -3.76%  '
$ws.Range('D28').Value = '''7.98'
$ws.Range('E28').Value = '  -2.42%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -3.11%  '
$ws.Range('E31').Value = '  -2.76%  '
$ws.Range('D32').Value = '''462.46'
$ws.Range('E32').Value = '  -0.29%  '
$ws.Range('D33').Value = '0.0₃0796'
$ws.Range('E33').Value = '  -3.84%  '
$ws.Range('D34').Value = '''1.65'
$ws.Range('E34').Value = '  +1.29%  '
$ws.Range('D35').Value = '''176.68'
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('D37').Value = '''0.399'
$ws.Range('E37').Value = '  -2.18%  '
$ws.Range('D38').Value = '''18.82'
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('E41').Value = '  -3.53%  '
$ws.Range('D42').Value = '''39.96'
$ws.Range('E42').Value = '  +1.00%  '
$ws.Range('D43').Value = '''157.97'
$ws.Range('E43').Value = '  +4.32%  '
$ws.Range('D45').Value = '''21.20'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('E46').Value = '  +2.80%  '
$ws.Range('D47').Value = '''0.0534'
$ws.Range('E47').Value = '  -2.73%  '
$ws.Range('E48').Value = '  -2.05%  '
$ws.Range('E49').Value = '  -1.66%  '
$ws.Range('D50').Value = '''17.99'
$ws.Range('E50').Value = '  -2.79%  '
$ws.Range('E51').Value = '  +0.15%  '
